# Refresh cryptos list snapshot (prices / 1h volume %) pulled from coinranking.com.
# Mirrors the per-run GitHub Actions data refresh for cryptos.xlsx.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "69.528.56"

# Row 3: Ethereum
$ws.Range("D3").Value = "3.692.09"
$ws.Range("E3").Value = "  +0.06%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "677.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.42%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.90%  "

# Row 7: USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8: XRP
$ws.Range("E8").Value = "  +0.28%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  +0.21%  "

# Row 10: Toncoin
$ws.Range("E10").Value = "  -0.91%  "

# Row 11: Cardano
$ws.Range("E11").Value = "  +0.75%  "

# Row 12: ShibaInu
$ws.Range("E12").Value = "  -0.24%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.312.51"

# Row 14: Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "

# Row 15: WrappedEther
$ws.Range("D15").Value = "3.682.29"
$ws.Range("E15").Value = "  -0.07%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "69.468.41"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17: TRON
$ws.Range("E17").Value = "  +2.30%  "

# Row 18: Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19: Polkadot
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.07%  "

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "

# Row 21: Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.07%  "

# Row 22: Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.650"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "

# Row 23: Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.20%  "

# Row 24: WrappedeETH
$ws.Range("D24").Value = "3.837.99"
$ws.Range("E24").Value = "  -0.01%  "

# Row 25: Dai
$ws.Range("E25").Value = "  -0.11%  "

# Row 26: PEPE
$ws.Range("E26").Value = "  +0.30%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("E27").Value = "  -1.21%  "

# Row 28: RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "

# Row 29: PancakeSwap
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "

# Row 30: Fetch.AI
$ws.Range("E30").Value = "  -0.28%  "

# Row 31: ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.33%  "

# Row 32: Binance-PegBSC-USD
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.77%  "

# Row 33: NEARProtocol
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.05%  "

# Row 34: EthereumClassic
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.01"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.69%  "

# Row 35: RenzoRestakedETH
$ws.Range("D35").Value = "3.681.38"
$ws.Range("E35").Value = "  +0.43%  "

# Row 36: Kaspa
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.162"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.97%  "

# Row 37: Aptos
$ws.Range("E37").Value = "  +3.42%  "

# Row 38: Filecoin
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.66%  "

# Row 40: Stacks
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.01%  "

# Row 42: Hedera
$ws.Range("E42").Value = "  -0.42%  "

# Row 43: Monero
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "168.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "

# Row 44: Mantle
$ws.Range("E44").Value = "  +0.01%  "

# Row 45: OKB
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.63%  "

# Row 46: dogwifhat
$ws.Range("E46").Value = "  -0.26%  "

# Row 47: FLOKI
$ws.Range("E47").Value = "  -0.34%  "

# Row 48: ONDO
$ws.Range("E48").Value = "  -0.95%  "

# Row 49: InjectiveProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "27.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.49%  "

# Row 50: SuiNetwork
$ws.Range("E50").Value = "  -2.47%  "

# Row 51: Cosmos
$ws.Range("E51").Value = "  +0.91%  "
